$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting rows 3-108 down to 4-109.
$ws.Rows("3:3").Insert()

# New row 3 repeats the same market/category metadata as the rest of the
# dataset, with fresh date/price figures for the latest week.
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 44860
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 100112040
$ws.Range("G3").Value = "Cilantro"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 600
$ws.Range("L3").Value = 700
$ws.Range("M3").Value = 650
$ws.Range("N3").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O3").Value = "Provincia de Diguillín"
$ws.Range("P3").Value = 650
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = "Hortaliza"
